# Update cryptos list with fresh values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = "'63.714.20"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3: D3, E3
$ws.Range("D3").Value = "'3.132.84"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4: E4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5: D5, E5
$ws.Range("D5").Value = "'585.97"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6: D6, E6
$ws.Range("D6").Value = "'146.31"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7: E7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8: D8, E8
$ws.Range("D8").Value = "'3.131.44"
$ws.Range("E8").Value = "  -0.48%  "

# Row 9: D9, E9
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -1.28%  "

# Row 10: E10
$ws.Range("E10").Value = "  +1.78%  "

# Row 11: D11, E11
$ws.Range("D11").Value = "'5.75"
$ws.Range("E11").Value = "  -0.89%  "

# Row 12: E12
$ws.Range("E12").Value = "  -3.15%  "

# Row 13: E13
$ws.Range("E13").Value = "  -2.81%  "

# Row 14: D14, E14
$ws.Range("D14").Value = "'36.83"
$ws.Range("E14").Value = "  +2.21%  "

# Row 15: E15
$ws.Range("E15").Value = "  -1.85%  "

# Row 16: D16, E16
$ws.Range("D16").Value = "'3.651.74"
$ws.Range("E16").Value = "  -0.87%  "

# Row 17: D17, E17
$ws.Range("D17").Value = "'63.548.06"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18: D18, E18
$ws.Range("D18").Value = "'3.133.97"
$ws.Range("E18").Value = "  -0.66%  "

# Row 19: D19, E19
$ws.Range("D19").Value = "'7.06"
$ws.Range("E19").Value = "  -1.97%  "

# Row 20: D20, E20
$ws.Range("D20").Value = "'463.59"
$ws.Range("E20").Value = "  -2.77%  "

# Row 21: D21, E21
$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22: D22, E22
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23: E23
$ws.Range("E23").Value = "  -3.42%  "

# Row 24: D24, E24
$ws.Range("D24").Value = "'12.93"
$ws.Range("E24").Value = "  -3.61%  "

# Row 25: D25, E25
$ws.Range("D25").Value = "'81.04"
$ws.Range("E25").Value = "  -2.24%  "

# Row 26: E26
$ws.Range("E26").Value = "  -1.87%  "

# Row 27: E27
$ws.Range("E27").Value = "  +0.02%  "

# Row 28: E28
$ws.Range("E28").Value = "  +5.61%  "

# Row 29: E29
$ws.Range("E29").Value = "  -1.34%  "

# Row 30: B30, C30, D30, E30
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  -1.26%  "

# Row 31: B31, C31, D31, E31
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.26%  "

# Row 32: D32, E32
$ws.Range("D32").Value = "'6.98"
$ws.Range("E32").Value = "  +0.82%  "

# Row 33: D33, E33
$ws.Range("D33").Value = "'26.92"
$ws.Range("E33").Value = "  -1.57%  "

# Row 34: E34
$ws.Range("E34").Value = "  -0.08%  "

# Row 35: D35, E35
$ws.Range("D35").Value = "'0.0₃0841"
$ws.Range("E35").Value = "  -5.31%  "

# Row 36: E36
$ws.Range("E36").Value = "  -1.84%  "

# Row 37: D37, E37
$ws.Range("D37").Value = "'2.29"
$ws.Range("E37").Value = "  -6.44%  "

# Row 38: E38
$ws.Range("E38").Value = "  -3.80%  "

# Row 39: E39
$ws.Range("E39").Value = "  -2.80%  "

# Row 40: D40, E40
$ws.Range("D40").Value = "'51.17"
$ws.Range("E40").Value = "  +0.55%  "

# Row 41: D41, E41
$ws.Range("D41").Value = "'438.72"
$ws.Range("E41").Value = "  -2.86%  "

# Row 42: D42, E42
$ws.Range("D42").Value = "'8.78"
$ws.Range("E42").Value = "  -0.11%  "

# Row 43: B43, C43, D43, E43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0371"
$ws.Range("E43").Value = "  -1.37%  "

# Row 44: B44, C44, D44, E44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'2.911.78"
$ws.Range("E44").Value = "  -1.12%  "

# Row 45: E45
$ws.Range("E45").Value = "  -2.65%  "

# Row 46: E46
$ws.Range("E46").Value = "  -4.02%  "

# Row 47: E47
$ws.Range("E47").Value = "  +3.74%  "

# Row 48: D48, E48
$ws.Range("D48").Value = "'126.68"
$ws.Range("E48").Value = "  +2.49%  "

# Row 49: E49
$ws.Range("E49").Value = "  -0.05%  "

# Row 51: E51
$ws.Range("E51").Value = "  -4.07%  "

